$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Iter3-JSON-All functions")
$ws.Range("A1").Value = "test"
